# Reproduce: Insert > Header & Footer, with "Slide Number" turned on and
# applied to all slides (Header / Footer / Date-and-time left off).
#
# This mirrors the author's commit: the PDO is "done" and a slide-number
# footer was switched on for every slide in the deck via Apply to All.

$p = $ppt.ActivePresentation

# Turn off header / footer / date-time on the slide master's header-footer
# defaults, and make sure slide number is on (mirrors the Header & Footer
# dialog "Apply to All" with only "Slide number" checked).
$masterHF = $p.SlideMaster.HeadersFooters
$masterHF.Header.Visible = $false
$masterHF.Footer.Visible = $false
$masterHF.DateAndTime.Visible = $false
$masterHF.SlideNumber.Visible = $true

# Apply the same settings to every individual slide so the slide-number
# placeholder actually gets instantiated (PowerPoint materialises a
# "Slide Number Placeholder" shape, bound to the layout's sldNum
# placeholder, the first time Apply to All turns slide numbers on).
$count = $p.Slides.Count
for ($i = 1; $i -le $count; $i++) {
    $slide = $p.Slides.Item($i)
    $hf = $slide.HeadersFooters
    $hf.Header.Visible = $false
    $hf.Footer.Visible = $false
    $hf.DateAndTime.Visible = $false
    $hf.SlideNumber.Visible = $true
}
